$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H13").Value = 85090
$ws_ALC.Range("I13").Value = 367201.34
$ws_ALC.Range("J13").Value = 456.6
$ws_ALC.Range("K13").Value = 367201.34
$ws_ALC.Range("L13").Value = 456.6
$ws_ALC.Range("M13").Value = -367032.34
$ws_ALC.Range("N13").Value = -794.6

$ws_ALC.Range("H17").Value = 358848.84
$ws_ALC.Range("J17").Value = 358848.84
$ws_ALC.Range("L17").Value = 1076546.52
$ws_ALC.Range("N17").Value = -1076882.52

$ws_ALC.Range("H57").Value = 43053
$ws_ALC.Range("J57").Value = 43053
$ws_ALC.Range("L57").Value = 129159
$ws_ALC.Range("N57").Value = -130157

$ws_ALC.Range("H62").Value = 4560.5454
$ws_ALC.Range("I62").Value = 3187.6667
$ws_ALC.Range("K62").Value = 3187.6667
$ws_ALC.Range("M62").Value = -2563.6667

$ws_ALC.Range("H65").Value = 4560.5454
$ws_ALC.Range("I65").Value = 3187.6667
$ws_ALC.Range("K65").Value = 15938.3335
$ws_ALC.Range("M65").Value = -12818.3335

$ws_ALC.Range("H132").Value = 2343.25
$ws_ALC.Range("I132").Value = 1896.1034
$ws_ALC.Range("K132").Value = 5688.3102
$ws_ALC.Range("M132").Value = -3158.3102

$ws_ALC.Range("H141").Value = 6297.5386
$ws_ALC.Range("I141").Value = 6047.3335
$ws_ALC.Range("J141").Value = 9300
$ws_ALC.Range("K141").Value = 18142.0005
$ws_ALC.Range("L141").Value = 27900
$ws_ALC.Range("M141").Value = -12962.0005
$ws_ALC.Range("N141").Value = -38260

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H9").Value = 1000
$ws_ARM.Range("I9").Value = 1000
$ws_ARM.Range("K9").Value = 1000
$ws_ARM.Range("M9").Value = -830

$ws_ARM.Range("H20").Value = 1000
$ws_ARM.Range("I20").Value = 1000
$ws_ARM.Range("K20").Value = 1000
$ws_ARM.Range("M20").Value = -730

$ws_ARM.Range("H32").Value = 20110.861
$ws_ARM.Range("I32").Value = 4753.355
$ws_ARM.Range("J32").Value = 337499.34
$ws_ARM.Range("K32").Value = 4753.355
$ws_ARM.Range("L32").Value = 337499.34
$ws_ARM.Range("M32").Value = -4466.355
$ws_ARM.Range("N32").Value = -338073.34

$ws_ARM.Range("H61").Value = 2363.2856
$ws_ARM.Range("I61").Value = 2290.4443
$ws_ARM.Range("K61").Value = 2290.4443
$ws_ARM.Range("M61").Value = -2078.4443

$ws_ARM.Range("H74").Value = 2703.111
$ws_ARM.Range("I74").Value = 2529.5334
$ws_ARM.Range("K74").Value = 2529.5334
$ws_ARM.Range("M74").Value = -1655.5334

$ws_ARM.Range("H76").Value = 25288
$ws_ARM.Range("J76").Value = 25288
$ws_ARM.Range("L76").Value = 25288
$ws_ARM.Range("N76").Value = -25964

$ws_ARM.Range("H77").Value = 2703.111
$ws_ARM.Range("I77").Value = 2529.5334
$ws_ARM.Range("K77").Value = 12647.667
$ws_ARM.Range("M77").Value = -8279.666999999999

$ws_ARM.Range("H79").Value = 25288
$ws_ARM.Range("J79").Value = 25288
$ws_ARM.Range("L79").Value = 25288
$ws_ARM.Range("N79").Value = -27628

$ws_ARM.Range("H110").Value = 2116.0527
$ws_ARM.Range("I110").Value = 2543.3572
$ws_ARM.Range("J110").Value = 919.6
$ws_ARM.Range("K110").Value = 2543.3572
$ws_ARM.Range("L110").Value = 919.6
$ws_ARM.Range("M110").Value = -498.3571999999999
$ws_ARM.Range("N110").Value = -5009.6

$ws_ARM.Range("H136").Value = 2363.2856
$ws_ARM.Range("I136").Value = 2290.4443
$ws_ARM.Range("K136").Value = 6871.3329
$ws_ARM.Range("M136").Value = -4321.3329

$ws_ARM.Range("H137").Value = 88973.16
$ws_ARM.Range("J137").Value = 88973.16
$ws_ARM.Range("L137").Value = 88973.16
$ws_ARM.Range("N137").Value = -99173.16

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 1654.341
$ws_BSM.Range("I94").Value = 1287.1212
$ws_BSM.Range("J94").Value = 2756
$ws_BSM.Range("K94").Value = 1287.1212
$ws_BSM.Range("L94").Value = 2756
$ws_BSM.Range("M94").Value = -836.1212
$ws_BSM.Range("N94").Value = -3658

$ws_BSM.Range("H105").Value = 1413.6154
$ws_BSM.Range("I105").Value = 1408.3
$ws_BSM.Range("K105").Value = 1408.3
$ws_BSM.Range("M105").Value = 338.7

$ws_BSM.Range("H134").Value = 835.63635
$ws_BSM.Range("I134").Value = 835.63635
$ws_BSM.Range("K134").Value = 2506.90905
$ws_BSM.Range("M134").Value = 28.09094999999979

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 142857700
$ws_CRP.Range("I7").Value = 200000700
$ws_CRP.Range("K7").Value = 200000700
$ws_CRP.Range("M7").Value = -200000587

$ws_CRP.Range("H21").Value = 2062.5
$ws_CRP.Range("J21").Value = 2062.5
$ws_CRP.Range("L21").Value = 2062.5
$ws_CRP.Range("N21").Value = -2532.5

$ws_CRP.Range("H22").Value = 0
$ws_CRP.Range("J22").Value = 0
$ws_CRP.Range("L22").Value = 0
$ws_CRP.Range("N22").ClearContents()

$ws_CRP.Range("H31").Value = 15872.77
$ws_CRP.Range("I31").Value = 2928.2104
$ws_CRP.Range("K31").Value = 2928.2104
$ws_CRP.Range("M31").Value = -2633.2104

$ws_CRP.Range("H34").Value = 15872.77
$ws_CRP.Range("I34").Value = 2928.2104
$ws_CRP.Range("K34").Value = 2928.2104
$ws_CRP.Range("M34").Value = -2726.2104

$ws_CRP.Range("H52").Value = 49975
$ws_CRP.Range("I52").Value = 49975
$ws_CRP.Range("K52").Value = 49975
$ws_CRP.Range("M52").Value = -49681

$ws_CRP.Range("H99").Value = 21083.076
$ws_CRP.Range("I99").Value = 28654.428
$ws_CRP.Range("K99").Value = 28654.428
$ws_CRP.Range("M99").Value = -27156.428

$ws_CRP.Range("H126").Value = 21083.076
$ws_CRP.Range("I126").Value = 28654.428
$ws_CRP.Range("K126").Value = 85963.284
$ws_CRP.Range("M126").Value = -83493.284

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 11180651
$ws_CUL.Range("I4").Value = 2624119.5
$ws_CUL.Range("K4").Value = 7872358.5
$ws_CUL.Range("M4").Value = -7872246.5

$ws_CUL.Range("H25").Value = 264446.75
$ws_CUL.Range("I25").Value = 400014.8
$ws_CUL.Range("K25").Value = 1200044.4
$ws_CUL.Range("M25").Value = -1199875.4

$ws_CUL.Range("H30").Value = 264446.75
$ws_CUL.Range("I30").Value = 400014.8
$ws_CUL.Range("K30").Value = 1200044.4
$ws_CUL.Range("M30").Value = -1199942.4

$ws_CUL.Range("H56").Value = 8110.696
$ws_CUL.Range("I56").Value = 8110.696
$ws_CUL.Range("K56").Value = 8110.696
$ws_CUL.Range("M56").Value = -7580.696

$ws_CUL.Range("H70").Value = 3163.6667
$ws_CUL.Range("I70").Value = 2497
$ws_CUL.Range("K70").Value = 7491
$ws_CUL.Range("M70").Value = -7176

$ws_CUL.Range("H73").Value = 3163.6667
$ws_CUL.Range("I73").Value = 2497
$ws_CUL.Range("K73").Value = 7491
$ws_CUL.Range("M73").Value = -6399

$ws_CUL.Range("H131").Value = 63983.125
$ws_CUL.Range("I131").Value = 101023
$ws_CUL.Range("K131").Value = 303069
$ws_CUL.Range("M131").Value = -298029

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 19301.555
$ws_GSM.Range("I97").Value = 24448.25
$ws_GSM.Range("J97").Value = 1288.125
$ws_GSM.Range("K97").Value = 24448.25
$ws_GSM.Range("L97").Value = 1288.125
$ws_GSM.Range("M97").Value = -23952.25
$ws_GSM.Range("N97").Value = -2280.125

$ws_GSM.Range("H122").Value = 1549.0588
$ws_GSM.Range("I122").Value = 1532.4615
$ws_GSM.Range("K122").Value = 4597.3845
$ws_GSM.Range("M122").Value = -2147.3845

$ws_GSM.Range("H124").Value = 46000
$ws_GSM.Range("J124").Value = 46000
$ws_GSM.Range("L124").Value = 46000
$ws_GSM.Range("N124").Value = -55820

$ws_GSM.Range("H132").Value = 5368.2383
$ws_GSM.Range("J132").Value = 5099.75
$ws_GSM.Range("L132").Value = 15299.25
$ws_GSM.Range("N132").Value = -20359.25

$ws_GSM.Range("H137").Value = 97313.8
$ws_GSM.Range("J137").Value = 97313.8
$ws_GSM.Range("L137").Value = 97313.8
$ws_GSM.Range("N137").Value = -107513.8

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H21").Value = 5999.5
$ws_LTW.Range("J21").Value = 5999.5
$ws_LTW.Range("L21").Value = 5999.5
$ws_LTW.Range("N21").Value = -6347.5

$ws_LTW.Range("H46").Value = 72686.836
$ws_LTW.Range("J46").Value = 3666.3333
$ws_LTW.Range("L46").Value = 3666.3333
$ws_LTW.Range("N46").Value = -4042.3333

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H20").Value = 23750
$ws_WVR.Range("I20").Value = 0
$ws_WVR.Range("J20").Value = 23750
$ws_WVR.Range("K20").Value = 0
$ws_WVR.Range("L20").Value = 23750
$ws_WVR.Range("M20").ClearContents()
$ws_WVR.Range("N20").Value = -24230

$ws_WVR.Range("H122").Value = 1154.1818
$ws_WVR.Range("I122").Value = 1125.5358
$ws_WVR.Range("K122").Value = 3376.6074
$ws_WVR.Range("M122").Value = -926.6074000000003

$ws_WVR.Range("H126").Value = 4502
$ws_WVR.Range("I126").Value = 4502
$ws_WVR.Range("K126").Value = 13506
$ws_WVR.Range("M126").Value = -11036

$ws_WVR.Range("H136").Value = 1630
$ws_WVR.Range("I136").Value = 1630
$ws_WVR.Range("K136").Value = 4890
$ws_WVR.Range("M136").Value = -2340
